{"js": "// Update the date title and the 25 \"two-digit \u00f7 one-digit\" practice\n// answers in the table, per the commit's regenerated problem set.\n\n// 1) Title paragraph: \"2024-06-10 Monday\" -> \"2024-06-11 Tuesday\"\nconst titleResults = context.document.body.search(\"2024-06-10 Monday\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length === 0) {\n  throw new Error(\"Title text '2024-06-10 Monday' not found\");\n}\ntitleResults.items[0].insertText(\"2024-06-11 Tuesday\", Word.InsertLocation.replace);\n\n// 2) Table cells: positional (row, col) updates, each paired with the\n// expected previous value so we never clobber the wrong cell even\n// though some new values collide with other cells' old values\n// elsewhere in the table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document\");\n}\nconst table = tables.items[0];\n\nconst cellUpdates = [[0, 0, \"89\u00f75=17, 4\", \"85\u00f79=9, 4\"], [0, 1, \"62\u00f76=10, 2\", \"83\u00f79=9, 2\"], [0, 2, \"55\u00f73=18, 1\", \"17\u00f72=8, 1\"], [0, 3, \"60\u00f72=30, 0\", \"28\u00f74=7, 0\"], [0, 4, \"19\u00f74=4, 3\", \"14\u00f79=1, 5\"], [4, 0, \"91\u00f74=22, 3\", \"58\u00f72=29, 0\"], [4, 1, \"74\u00f73=24, 2\", \"10\u00f75=2, 0\"], [4, 2, \"30\u00f77=4, 2\", \"41\u00f77=5, 6\"], [4, 3, \"56\u00f78=7, 0\", \"41\u00f77=5, 6\"], [4, 4, \"42\u00f79=4, 6\", \"28\u00f74=7, 0\"], [8, 0, \"55\u00f72=27, 1\", \"30\u00f72=15, 0\"], [8, 1, \"87\u00f74=21, 3\", \"72\u00f73=24, 0\"], [8, 2, \"63\u00f76=10, 3\", \"94\u00f77=13, 3\"], [8, 3, \"18\u00f75=3, 3\", \"39\u00f76=6, 3\"], [8, 4, \"31\u00f73=10, 1\", \"50\u00f73=16, 2\"], [12, 0, \"31\u00f78=3, 7\", \"66\u00f72=33, 0\"], [12, 1, \"76\u00f79=8, 4\", \"12\u00f73=4, 0\"], [12, 2, \"69\u00f76=11, 3\", \"73\u00f73=24, 1\"], [12, 3, \"96\u00f79=10, 6\", \"78\u00f78=9, 6\"], [12, 4, \"83\u00f79=9, 2\", \"54\u00f78=6, 6\"], [16, 0, \"61\u00f74=15, 1\", \"36\u00f75=7, 1\"], [16, 1, \"35\u00f79=3, 8\", \"78\u00f79=8, 6\"], [16, 2, \"24\u00f79=2, 6\", \"96\u00f77=13, 5\"], [16, 3, \"89\u00f79=9, 8\", \"84\u00f74=21, 0\"], [16, 4, \"53\u00f76=8, 5\", \"40\u00f75=8, 0\"]];\n\nconst cells = cellUpdates.map(([row, col]) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\ncellUpdates.forEach(([row, col, oldText, newText], i) => {\n  const cell = cells[i];\n  if (cell.value !== oldText) {\n    throw new Error(\n      `Cell (${row},${col}) expected \"${oldText}\" but found \"${cell.value}\"`\n    );\n  }\n  cell.value = newText;\n});\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 \"two-digit \u00f7 one-digit\" practice\n# answers in the table, per the commit's regenerated problem set.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2024-06-10 Monday\" -> \"2024-06-11 Tuesday\"\n$titleRange = $d.Content.Duplicate\n$titleRange.Find.ClearFormatting()\n$titleRange.Find.Text = \"2024-06-10 Monday\"\n$titleRange.Find.Replacement.ClearFormatting()\n$titleRange.Find.Replacement.Text = \"2024-06-11 Tuesday\"\n$titleFound = $titleRange.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\nif (-not $titleFound) {\n    throw \"Title text '2024-06-10 Monday' not found\"\n}\n\n# 2) Table cells: positional (row, col) updates (1-based, as in the Word\n# object model), each paired with the expected previous value so we\n# never clobber the wrong cell even though some new values collide\n# with other cells' old values elsewhere in the table.\n$cellUpdates = @(\n    @{ Row = 1; Col = 1; Old = \"89\u00f75=17, 4\"; New = \"85\u00f79=9, 4\" },\n    @{ Row = 1; Col = 2; Old = \"62\u00f76=10, 2\"; New = \"83\u00f79=9, 2\" },\n    @{ Row = 1; Col = 3; Old = \"55\u00f73=18, 1\"; New = \"17\u00f72=8, 1\" },\n    @{ Row = 1; Col = 4; Old = \"60\u00f72=30, 0\"; New = \"28\u00f74=7, 0\" },\n    @{ Row = 1; Col = 5; Old = \"19\u00f74=4, 3\"; New = \"14\u00f79=1, 5\" },\n    @{ Row = 5; Col = 1; Old = \"91\u00f74=22, 3\"; New = \"58\u00f72=29, 0\" },\n    @{ Row = 5; Col = 2; Old = \"74\u00f73=24, 2\"; New = \"10\u00f75=2, 0\" },\n    @{ Row = 5; Col = 3; Old = \"30\u00f77=4, 2\"; New = \"41\u00f77=5, 6\" },\n    @{ Row = 5; Col = 4; Old = \"56\u00f78=7, 0\"; New = \"41\u00f77=5, 6\" },\n    @{ Row = 5; Col = 5; Old = \"42\u00f79=4, 6\"; New = \"28\u00f74=7, 0\" },\n    @{ Row = 9; Col = 1; Old = \"55\u00f72=27, 1\"; New = \"30\u00f72=15, 0\" },\n    @{ Row = 9; Col = 2; Old = \"87\u00f74=21, 3\"; New = \"72\u00f73=24, 0\" },\n    @{ Row = 9; Col = 3; Old = \"63\u00f76=10, 3\"; New = \"94\u00f77=13, 3\" },\n    @{ Row = 9; Col = 4; Old = \"18\u00f75=3, 3\"; New = \"39\u00f76=6, 3\" },\n    @{ Row = 9; Col = 5; Old = \"31\u00f73=10, 1\"; New = \"50\u00f73=16, 2\" },\n    @{ Row = 13; Col = 1; Old = \"31\u00f78=3, 7\"; New = \"66\u00f72=33, 0\" },\n    @{ Row = 13; Col = 2; Old = \"76\u00f79=8, 4\"; New = \"12\u00f73=4, 0\" },\n    @{ Row = 13; Col = 3; Old = \"69\u00f76=11, 3\"; New = \"73\u00f73=24, 1\" },\n    @{ Row = 13; Col = 4; Old = \"96\u00f79=10, 6\"; New = \"78\u00f78=9, 6\" },\n    @{ Row = 13; Col = 5; Old = \"83\u00f79=9, 2\"; New = \"54\u00f78=6, 6\" },\n    @{ Row = 17; Col = 1; Old = \"61\u00f74=15, 1\"; New = \"36\u00f75=7, 1\" },\n    @{ Row = 17; Col = 2; Old = \"35\u00f79=3, 8\"; New = \"78\u00f79=8, 6\" },\n    @{ Row = 17; Col = 3; Old = \"24\u00f79=2, 6\"; New = \"96\u00f77=13, 5\" },\n    @{ Row = 17; Col = 4; Old = \"89\u00f79=9, 8\"; New = \"84\u00f74=21, 0\" },\n    @{ Row = 17; Col = 5; Old = \"53\u00f76=8, 5\"; New = \"40\u00f75=8, 0\" }\n)\n\n$table = $d.Tables.Item(1)\n\nforeach ($update in $cellUpdates) {\n    $cell = $table.Cell($update.Row, $update.Col)\n    $current = $cell.Range.Text.TrimEnd([char]7, [char]13)\n    if ($current -ne $update.Old) {\n        throw \"Cell ($($update.Row),$($update.Col)) expected '$($update.Old)' but found '$current'\"\n    }\n    $cell.Range.Text = $update.New\n}\n"}
